$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 151, shifting existing rows 151-244 down to 152-245
$ws.Rows("151").Insert()

# Populate the newly inserted row 151 with the new record
$ws.Range("A151").Value = 11
$ws.Range("B151").Value = "Vega Monumental Concepción"
$ws.Range("C151").Value = "Bíobío"
$ws.Range("D151").Value = 44762
$ws.Range("E151").Value = 8
$ws.Range("F151").Value = 100114013
$ws.Range("G151").Value = "Zanahoria"
$ws.Range("H151").Value = "Sin especificar"
$ws.Range("I151").Value = "Primera"
$ws.Range("J151").Value = 280
$ws.Range("K151").Value = 9000
$ws.Range("L151").Value = 10000
$ws.Range("M151").Value = 9643
$ws.Range("N151").Value = "$/saco 20 kilos"
$ws.Range("O151").Value = "Región de La Araucanía"
$ws.Range("P151").Value = 482
$ws.Range("Q151").Value = 20
$ws.Range("R151").Value = "Hortaliza"
